$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028185895818088
$ws.Cells.Item(2, 4).Value = 1.033314673539573
$ws.Cells.Item(2, 5).Value = 1.028181577463523
$ws.Cells.Item(2, 6).Value = 1.026732192368816
$ws.Cells.Item(2, 9).Value = 1.035243628196483
$ws.Cells.Item(2, 10).Value = 1.03333961427441
$ws.Cells.Item(2, 11).Value = 1.036117378245251
$ws.Cells.Item(2, 12).Value = 1.030999140757189
$ws.Cells.Item(2, 13).Value = 1.029553980392997
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029130216393325
$ws.Cells.Item(3, 4).Value = 1.034040460955386
$ws.Cells.Item(3, 5).Value = 1.028983265508087
$ws.Cells.Item(3, 6).Value = 1.028309723001516
$ws.Cells.Item(3, 9).Value = 1.035506269172334
$ws.Cells.Item(3, 10).Value = 1.033924463486449
$ws.Cells.Item(3, 11).Value = 1.036652342380609
$ws.Cells.Item(3, 12).Value = 1.031608727094165
$ws.Cells.Item(3, 13).Value = 1.030937004167044
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.029741101974341
$ws.Cells.Item(4, 4).Value = 1.034509842357399
$ws.Cells.Item(4, 5).Value = 1.029502257531264
$ws.Cells.Item(4, 6).Value = 1.029330308996242
$ws.Cells.Item(4, 9).Value = 1.035674727158373
$ws.Cells.Item(4, 10).Value = 1.034302147432107
$ws.Cells.Item(4, 11).Value = 1.036997581361546
$ws.Cells.Item(4, 12).Value = 1.032002766181187
$ws.Cells.Item(4, 13).Value = 1.031831258504808
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.029997881914666
$ws.Cells.Item(5, 4).Value = 1.034707109708846
$ws.Cells.Item(5, 5).Value = 1.029720500419358
$ws.Cells.Item(5, 6).Value = 1.029759325292152
$ws.Cells.Item(5, 9).Value = 1.035745190567728
$ws.Cells.Item(5, 10).Value = 1.034460745724082
$ws.Cells.Item(5, 11).Value = 1.037142499865464
$ws.Cells.Item(5, 12).Value = 1.032168323436851
$ws.Cells.Item(5, 13).Value = 1.032207051008634
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030040994272116
$ws.Cells.Item(6, 4).Value = 1.034740228201565
$ws.Cells.Item(6, 5).Value = 1.029757147790179
$ws.Cells.Item(6, 6).Value = 1.029831356943115
$ws.Cells.Item(6, 9).Value = 1.035757000788077
$ws.Cells.Item(6, 10).Value = 1.034487364521462
$ws.Cells.Item(6, 11).Value = 1.037166819419551
$ws.Cells.Item(6, 12).Value = 1.032196115551998
$ws.Cells.Item(6, 13).Value = 1.032270139444749
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.02974453322301
$ws.Cells.Item(7, 4).Value = 1.03451247849195
$ws.Cells.Item(7, 5).Value = 1.029505173473401
$ws.Cells.Item(7, 6).Value = 1.029336041675722
$ws.Cells.Item(7, 9).Value = 1.035675670094066
$ws.Cells.Item(7, 10).Value = 1.03430426733752
$ws.Cells.Item(7, 11).Value = 1.036999518634075
$ws.Cells.Item(7, 12).Value = 1.032004978745682
$ws.Cells.Item(7, 13).Value = 1.031836280454168
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028505065213616
$ws.Cells.Item(8, 4).Value = 1.033560008674356
$ws.Cells.Item(8, 5).Value = 1.028452460518837
$ws.Cells.Item(8, 6).Value = 1.027265367172565
$ws.Cells.Item(8, 9).Value = 1.035332697087423
$ws.Cells.Item(8, 10).Value = 1.033537422718714
$ws.Cells.Item(8, 11).Value = 1.036298362049545
$ws.Cells.Item(8, 12).Value = 1.031205236943965
$ws.Cells.Item(8, 13).Value = 1.030021518092769
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026319779978969
$ws.Cells.Item(9, 4).Value = 1.031879723389283
$ws.Cells.Item(9, 5).Value = 1.026599344785221
$ws.Cells.Item(9, 6).Value = 1.0236149366482
$ws.Cells.Item(9, 9).Value = 1.034716937151585
$ws.Cells.Item(9, 10).Value = 1.032180373729709
$ws.Cells.Item(9, 11).Value = 1.035055800882851
$ws.Cells.Item(9, 12).Value = 1.029792896781966
$ws.Cells.Item(9, 13).Value = 1.026818452567219
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024862106045968
$ws.Cells.Item(10, 4).Value = 1.030758265516502
$ws.Cells.Item(10, 5).Value = 1.025365225911733
$ws.Cells.Item(10, 6).Value = 1.021179871528655
$ws.Cells.Item(10, 9).Value = 1.034298766789389
$ws.Cells.Item(10, 10).Value = 1.031271782480322
$ws.Cells.Item(10, 11).Value = 1.034222697327133
$ws.Cells.Item(10, 12).Value = 1.02884925661699
$ws.Cells.Item(10, 13).Value = 1.024679269998502
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024230717670703
$ws.Cells.Item(11, 4).Value = 1.030272364074313
$ws.Cells.Item(11, 5).Value = 1.024831147435686
$ws.Cells.Item(11, 6).Value = 1.02012504188586
$ws.Cells.Item(11, 9).Value = 1.034115875981357
$ws.Cells.Item(11, 10).Value = 1.030877425873769
$ws.Cells.Item(11, 11).Value = 1.033860831616493
$ws.Cells.Item(11, 12).Value = 1.028440155700229
$ws.Cells.Item(11, 13).Value = 1.023752011775771
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.023996160373102
$ws.Cells.Item(12, 4).Value = 1.030091833289079
$ws.Cells.Item(12, 5).Value = 1.024632812542575
$ws.Cells.Item(12, 6).Value = 1.019733159843322
$ws.Cells.Item(12, 9).Value = 1.034047668579912
$ws.Cells.Item(12, 10).Value = 1.030730804068822
$ws.Cells.Item(12, 11).Value = 1.033726249178814
$ws.Cells.Item(12, 12).Value = 1.028288122337454
$ws.Cells.Item(12, 13).Value = 1.023407434091709
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024046475146766
$ws.Cells.Item(13, 4).Value = 1.030130559805678
$ws.Cells.Item(13, 5).Value = 1.024675353984096
$ws.Cells.Item(13, 6).Value = 1.01981722320801
$ws.Cells.Item(13, 9).Value = 1.034062311682797
$ws.Cells.Item(13, 10).Value = 1.030762261294511
$ws.Cells.Item(13, 11).Value = 1.033755125237186
$ws.Cells.Item(13, 12).Value = 1.028320737409665
$ws.Cells.Item(13, 13).Value = 1.023481354176991
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024211329738275
$ws.Cells.Item(14, 4).Value = 1.030257442267768
$ws.Cells.Item(14, 5).Value = 1.024814752080172
$ws.Cells.Item(14, 6).Value = 1.020092650288217
$ws.Cells.Item(14, 9).Value = 1.034110243518676
$ws.Cells.Item(14, 10).Value = 1.030865308932007
$ws.Cells.Item(14, 11).Value = 1.033849710446977
$ws.Cells.Item(14, 12).Value = 1.02842759010741
$ws.Cells.Item(14, 13).Value = 1.023723532025188
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024312897895187
$ws.Cells.Item(15, 4).Value = 1.030335612734222
$ws.Cells.Item(15, 5).Value = 1.02490064589718
$ws.Cells.Item(15, 6).Value = 1.020262340347861
$ws.Cells.Item(15, 9).Value = 1.034139739644337
$ws.Cells.Item(15, 10).Value = 1.030928781392465
$ws.Cells.Item(15, 11).Value = 1.033907965063816
$ws.Cells.Item(15, 12).Value = 1.028493415628646
$ws.Cells.Item(15, 13).Value = 1.023872725384397
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.024904004822817
$ws.Cells.Item(16, 4).Value = 1.030790506827391
$ws.Cells.Item(16, 5).Value = 1.025400677335895
$ws.Cells.Item(16, 6).Value = 1.021249867460341
$ws.Cells.Item(16, 9).Value = 1.034310866287269
$ws.Cells.Item(16, 10).Value = 1.031297934989692
$ws.Cells.Item(16, 11).Value = 1.034246689416344
$ws.Cells.Item(16, 12).Value = 1.028876396780086
$ws.Cells.Item(16, 13).Value = 1.024740787975994
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025274735022631
$ws.Cells.Item(17, 4).Value = 1.031075768908843
$ws.Cells.Item(17, 5).Value = 1.025714415072168
$ws.Cells.Item(17, 6).Value = 1.021869197193426
$ws.Cells.Item(17, 9).Value = 1.034417722073953
$ws.Cells.Item(17, 10).Value = 1.031529245867326
$ws.Cells.Item(17, 11).Value = 1.034458860570038
$ws.Cells.Item(17, 12).Value = 1.029116497156256
$ws.Cells.Item(17, 13).Value = 1.025285034801428
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025490955846198
$ws.Cells.Item(18, 4).Value = 1.031242128232066
$ws.Cells.Item(18, 5).Value = 1.025897442228872
$ws.Cells.Item(18, 6).Value = 1.022230400673473
$ws.Cells.Item(18, 9).Value = 1.034479873624972
$ws.Cells.Item(18, 10).Value = 1.031664075826003
$ws.Cells.Item(18, 11).Value = 1.034582507747849
$ws.Cells.Item(18, 12).Value = 1.0292564955577
$ws.Cells.Item(18, 13).Value = 1.0256023907859
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025564678214113
$ws.Cells.Item(19, 4).Value = 1.031298847501257
$ws.Cells.Item(19, 5).Value = 1.025959854716748
$ws.Cells.Item(19, 6).Value = 1.022353554961806
$ws.Cells.Item(19, 9).Value = 1.034501035913339
$ws.Cells.Item(19, 10).Value = 1.031710034144033
$ws.Cells.Item(19, 11).Value = 1.034624649813668
$ws.Cells.Item(19, 12).Value = 1.029304223243681
$ws.Cells.Item(19, 13).Value = 1.025710585209229
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02523496123599
$ws.Cells.Item(20, 4).Value = 1.031045166008724
$ws.Cells.Item(20, 5).Value = 1.025680750956715
$ws.Cells.Item(20, 6).Value = 1.021802753270503
$ws.Cells.Item(20, 9).Value = 1.034406275625563
$ws.Cells.Item(20, 10).Value = 1.031504437686521
$ws.Cells.Item(20, 11).Value = 1.03443610785424
$ws.Cells.Item(20, 12).Value = 1.029090741633043
$ws.Cells.Item(20, 13).Value = 1.025226652029547
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024162785045563
$ws.Cells.Item(21, 4).Value = 1.030220079796126
$ws.Cells.Item(21, 5).Value = 1.024773701554925
$ws.Cells.Item(21, 6).Value = 1.020011545901315
$ws.Cells.Item(21, 9).Value = 1.03409613634182
$ws.Cells.Item(21, 10).Value = 1.030834967844609
$ws.Cells.Item(21, 11).Value = 1.033821862137905
$ws.Cells.Item(21, 12).Value = 1.028396126724806
$ws.Cells.Item(21, 13).Value = 1.023652220946443
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023488482521965
$ws.Cells.Item(22, 4).Value = 1.029701052997207
$ws.Cells.Item(22, 5).Value = 1.024203667649822
$ws.Cells.Item(22, 6).Value = 1.018884926611479
$ws.Cells.Item(22, 9).Value = 1.03389955635958
$ws.Cells.Item(22, 10).Value = 1.030413234216387
$ws.Cells.Item(22, 11).Value = 1.033434681284172
$ws.Cells.Item(22, 12).Value = 1.027958960356571
$ws.Cells.Item(22, 13).Value = 1.022661427328571
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023845960578084
$ws.Cells.Item(23, 4).Value = 1.029976223748947
$ws.Cells.Item(23, 5).Value = 1.024505828351661
$ws.Cells.Item(23, 6).Value = 1.019482210498729
$ws.Cells.Item(23, 9).Value = 1.034003917233916
$ws.Cells.Item(23, 10).Value = 1.030636880232874
$ws.Cells.Item(23, 11).Value = 1.033640026152118
$ws.Cells.Item(23, 12).Value = 1.02819075175648
$ws.Cells.Item(23, 13).Value = 1.02318675171167
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.02525293335914
$ws.Cells.Item(24, 4).Value = 1.031058994232184
$ws.Cells.Item(24, 5).Value = 1.02569596223096
$ws.Cells.Item(24, 6).Value = 1.021832776544583
$ws.Cells.Item(24, 9).Value = 1.034411448325417
$ws.Cells.Item(24, 10).Value = 1.031515647712729
$ws.Cells.Item(24, 11).Value = 1.034446389162739
$ws.Cells.Item(24, 12).Value = 1.029102379593777
$ws.Cells.Item(24, 13).Value = 1.025253032980122
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.026884871056712
$ws.Cells.Item(25, 4).Value = 1.032314342108265
$ws.Cells.Item(25, 5).Value = 1.027078193018232
$ws.Cells.Item(25, 6).Value = 1.024558889811961
$ws.Cells.Item(25, 9).Value = 1.034877476583387
$ws.Cells.Item(25, 10).Value = 1.032531888495359
$ws.Cells.Item(25, 11).Value = 1.035377866268381
$ws.Cells.Item(25, 12).Value = 1.030158387123597
$ws.Cells.Item(25, 13).Value = 1.027647170896794
